$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("Users")
$wsBirds = $wb.Worksheets.Item("Birds")
$wsCages = $wb.Worksheets.Item("Cages")

# --- Users sheet: row 6 text correction (" " -> Hebrew "Shin") ---
$wsUsers.Range("A6").Value = "ש"
$wsUsers.Range("B6").Value = "ש"

# --- Users sheet: new login rows, interleaved with new cage rows so the
#     shared-string table ends up in the same order as the target file ---

# Row 7: kakape98 / kaka1998!  (C7 also updates from 3 -> 334)
$wsUsers.Range("A7").Value = "kakape98"
$wsUsers.Range("B7").Value = "kaka1998!"
$wsUsers.Range("C7").Value = 334

# Row 8: avivaa98 / avivaa98!
$wsUsers.Range("A8").Value = "avivaa98"
$wsUsers.Range("B8").Value = "avivaa98!"
$wsUsers.Range("C8").Value = 106

# Cages row 9: a9
$wsCages.Range("A9").Value = "a9"

# Row 9: aaaaaa11 / aaaaaa1!
$wsUsers.Range("A9").Value = "aaaaaa11"
$wsUsers.Range("B9").Value = "aaaaaa1!"
$wsUsers.Range("C9").Value = 102

# Cages row 10: a10
$wsCages.Range("A10").Value = "a10"

# Row 10: bbbbbb9 / bbbbbb1!
$wsUsers.Range("A10").Value = "bbbbbb9"
$wsUsers.Range("B10").Value = "bbbbbb1!"
$wsUsers.Range("C10").Value = 1113

# --- Cages sheet: finish filling the new rows (reuses existing strings) ---
$wsCages.Range("B9").Value = 6
$wsCages.Range("C9").Value = 6
$wsCages.Range("D9").Value = 6
$wsCages.Range("E9").Value = "wood"

$wsCages.Range("B10").Value = 9
$wsCages.Range("C10").Value = 9
$wsCages.Range("D10").Value = 9
$wsCages.Range("E10").Value = "wood"

# --- Birds sheet: three new birds, each linked to a cage ---

# Row 7 (bird 6): American Gouldian / Notrh America / Male, cage a1
$wsBirds.Range("A7").Value = 6
$wsBirds.Range("B7").Value = "American Gouldian"
$wsBirds.Range("C7").Value = "Notrh America"
$wsBirds.Range("D7").Value = "Male"
$wsBirds.Range("G2").Copy($wsBirds.Range("G7"))
$wsBirds.Range("G7").Value = 45049
$wsBirds.Range("H7").Value = "a1"
$wsBirds.Range("I7").Value = 106

# Row 8 (bird 7): European Gouldian / East Europe / Female, cage a1
$wsBirds.Range("A8").Value = 7
$wsBirds.Range("B8").Value = "European Gouldian"
$wsBirds.Range("C8").Value = "East Europe"
$wsBirds.Range("D8").Value = "Female"
$wsBirds.Range("G2").Copy($wsBirds.Range("G8"))
$wsBirds.Range("G8").Value = 45049
$wsBirds.Range("H8").Value = "a1"
$wsBirds.Range("I8").Value = 102

# Row 9 (bird 8): American Gouldian / Central America (new) / Male, cage a3
$wsBirds.Range("A9").Value = 8
$wsBirds.Range("B9").Value = "American Gouldian"
$wsBirds.Range("C9").Value = "Central America"
$wsBirds.Range("D9").Value = "Male"
$wsBirds.Range("G2").Copy($wsBirds.Range("G9"))
$wsBirds.Range("G9").Value = 45049
$wsBirds.Range("H9").Value = "a3"
$wsBirds.Range("I9").Value = 1
